$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the two new rows (16 and 17) inherit the same formatting as the
# other data rows (bold/bordered/centered style on column A, like A2:A15).
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A17").PasteSpecial(-4122)

# Rows 8 and 9 become the new "line7"/"line8" entries, and what used to be
# extr1..extr8 (rows 8-15) shift down to rows 10-17 with updated values.
$data = @(
    @{ Row = 8;  A = 6;  B = "line7"; C = 14; D = 11; E = $true  },
    @{ Row = 9;  A = 7;  B = "line8"; C = 16; D = 9;  E = $true  },
    @{ Row = 10; A = 8;  B = "extr1"; C = 5;  D = 12; E = $true  },
    @{ Row = 11; A = 9;  B = "extr2"; C = 5;  D = 9;  E = $true  },
    @{ Row = 12; A = 10; B = "extr3"; C = 10; D = 11; E = $true  },
    @{ Row = 13; A = 11; B = "extr4"; C = 7;  D = 8;  E = $false },
    @{ Row = 14; A = 12; B = "extr5"; C = 9;  D = 11; E = $true  },
    @{ Row = 15; A = 13; B = "extr6"; C = 7;  D = 11; E = $true  },
    @{ Row = 16; A = 14; B = "extr7"; C = 5;  D = 7;  E = $false },
    @{ Row = 17; A = 15; B = "extr8"; C = 8;  D = 5;  E = $true  }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
}
